$d = $word.ActiveDocument

# The document has a single section whose headers/footers each contain
# one inline logo picture. We need to rename the picture's "name"
# attribute (stored on both wp:docPr and pic:cNvPr) for each logo:
#   - the BTEC logo pictures (currently "image1.jpg") become "image2.jpg"
#   - the Pearson logo pictures (currently "image2.png") become "image1.png"

foreach ($sec in $d.Sections) {
    for ($hi = 1; $hi -le 3; $hi++) {
        $hf = $sec.Headers.Item($hi)
        if ($hf.Exists) {
            $count = $hf.Range.InlineShapes.Count
            for ($j = 1; $j -le $count; $j++) {
                $shp = $hf.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
    for ($fi = 1; $fi -le 3; $fi++) {
        $ft = $sec.Footers.Item($fi)
        if ($ft.Exists) {
            $count = $ft.Range.InlineShapes.Count
            for ($j = 1; $j -le $count; $j++) {
                $shp = $ft.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
